$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 1.433944
$ws.Cells.Item(2, 8).Value = 4.301832
$ws.Cells.Item(2, 9).Value = 0.003882998715548277
$ws.Cells.Item(2, 10).Value = 0.003886188105009087
$ws.Cells.Item(2, 13).Value = 1.986769333333333
$ws.Cells.Item(2, 14).Value = 5.960307999999999
$ws.Cells.Item(2, 15).Value = 0.4188172894723593
$ws.Cells.Item(2, 16).Value = 0.4376286247783243
$ws.Cells.Item(2, 17).Value = 2.848915964917333
$ws.Cells.Item(2, 18).Value = 25.640243684256
$ws.Cells.Item(2, 19).Value = 0.001626266997070582
$ws.Cells.Item(2, 20).Value = 0.001700707156025009

$ws.Cells.Item(3, 7).Value = 1.433944
$ws.Cells.Item(3, 8).Value = 4.301832
$ws.Cells.Item(3, 9).Value = 0.003882998715548277
$ws.Cells.Item(3, 10).Value = 0.003886188105009087
$ws.Cells.Item(3, 15).Value = 0.2760769767479527
$ws.Cells.Item(3, 16).Value = 0.2884770774849725
$ws.Cells.Item(3, 17).Value = 1.877955199973333
$ws.Cells.Item(3, 18).Value = 16.90159679976
$ws.Cells.Item(3, 19).Value = 0.001072006546104752
$ws.Cells.Item(3, 20).Value = 0.001121076187089885

$ws.Cells.Item(4, 7).Value = 1.433944
$ws.Cells.Item(4, 8).Value = 4.301832
$ws.Cells.Item(4, 9).Value = 0.003882998715548277
$ws.Cells.Item(4, 10).Value = 0.003886188105009087
$ws.Cells.Item(4, 13).Value = 0.3499836666666667
$ws.Cells.Item(4, 14).Value = 1.049951
$ws.Cells.Item(4, 15).Value = 0.07377766919071853
$ws.Cells.Item(4, 16).Value = 0.07709142081493547
$ws.Cells.Item(4, 17).Value = 0.5018569789146667
$ws.Cells.Item(4, 18).Value = 4.516712810232001
$ws.Cells.Item(4, 19).Value = 0.0002864785947037058
$ws.Cells.Item(4, 20).Value = 0.0002995917625692522

$ws.Cells.Item(5, 7).Value = 1.433944
$ws.Cells.Item(5, 8).Value = 4.301832
$ws.Cells.Item(5, 9).Value = 0.003882998715548277
$ws.Cells.Item(5, 10).Value = 0.003886188105009087
$ws.Cells.Item(5, 13).Value = 0.6117275
$ws.Cells.Item(5, 14).Value = 1.223455
$ws.Cells.Item(5, 15).Value = 0.1289541010862372
$ws.Cells.Item(5, 16).Value = 0.08983074853315715
$ws.Cells.Item(5, 17).Value = 0.87718297826
$ws.Cells.Item(5, 18).Value = 5.26309786956
$ws.Cells.Item(5, 19).Value = 0.0005007286088825416
$ws.Cells.Item(5, 20).Value = 0.0003490991864136178

$ws.Cells.Item(6, 7).Value = 1.433944
$ws.Cells.Item(6, 8).Value = 4.301832
$ws.Cells.Item(6, 9).Value = 0.003882998715548277
$ws.Cells.Item(6, 10).Value = 0.003886188105009087
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.4856376666666667
$ws.Cells.Item(6, 14).Value = 1.456913
$ws.Cells.Item(6, 15).Value = 0.1023739635027323
$ws.Cells.Item(6, 16).Value = 0.1069721283886106
$ws.Cells.Item(6, 17).Value = 0.6963772182906668
$ws.Cells.Item(6, 18).Value = 6.267394964616001
$ws.Cells.Item(6, 19).Value = 0.0003975179687866958
$ws.Cells.Item(6, 20).Value = 0.0004157138129113234

$ws.Cells.Item(7, 9).Value = 0.7877262822264709
$ws.Cells.Item(7, 10).Value = 0.7883732991550308
$ws.Cells.Item(7, 13).Value = 1.986769333333333
$ws.Cells.Item(7, 14).Value = 5.960307999999999
$ws.Cells.Item(7, 15).Value = 0.4188172894723593
$ws.Cells.Item(7, 16).Value = 0.4376286247783243
$ws.Cells.Item(7, 17).Value = 577.9466195633534
$ws.Cells.Item(7, 18).Value = 5201.51957607018
$ws.Cells.Item(7, 19).Value = 0.3299133863682293
$ws.Cells.Item(7, 20).Value = 0.3450147227211666

$ws.Cells.Item(8, 9).Value = 0.7877262822264709
$ws.Cells.Item(8, 10).Value = 0.7883732991550308
$ws.Cells.Item(8, 15).Value = 0.2760769767479527
$ws.Cells.Item(8, 16).Value = 0.2884770774849725
$ws.Cells.Item(8, 17).Value = 380.9722269387834
$ws.Cells.Item(8, 18).Value = 3428.75004244905
$ws.Cells.Item(8, 19).Value = 0.2174730905019887
$ws.Cells.Item(8, 20).Value = 0.2274276253074293

$ws.Cells.Item(9, 9).Value = 0.7877262822264709
$ws.Cells.Item(9, 10).Value = 0.7883732991550308
$ws.Cells.Item(9, 13).Value = 0.3499836666666667
$ws.Cells.Item(9, 14).Value = 1.049951
$ws.Cells.Item(9, 15).Value = 0.07377766919071853
$ws.Cells.Item(9, 16).Value = 0.07709142081493547
$ws.Cells.Item(9, 17).Value = 101.8094419209817
$ws.Cells.Item(9, 18).Value = 916.2849772888352
$ws.Cells.Item(9, 19).Value = 0.05811660906293915
$ws.Cells.Item(9, 20).Value = 0.0607768177644195

$ws.Cells.Item(10, 9).Value = 0.7877262822264709
$ws.Cells.Item(10, 10).Value = 0.7883732991550308
$ws.Cells.Item(10, 13).Value = 0.6117275
$ws.Cells.Item(10, 14).Value = 1.223455
$ws.Cells.Item(10, 15).Value = 0.1289541010862372
$ws.Cells.Item(10, 16).Value = 0.08983074853315715
$ws.Cells.Item(10, 17).Value = 177.9501197181125
$ws.Cells.Item(10, 18).Value = 1067.700718308675
$ws.Cells.Item(10, 19).Value = 0.1015805346265181
$ws.Cells.Item(10, 20).Value = 0.07082016358665105

$ws.Cells.Item(11, 9).Value = 0.7877262822264709
$ws.Cells.Item(11, 10).Value = 0.7883732991550308
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.4856376666666667
$ws.Cells.Item(11, 14).Value = 1.456913
$ws.Cells.Item(11, 15).Value = 0.1023739635027323
$ws.Cells.Item(11, 16).Value = 0.1069721283886106
$ws.Cells.Item(11, 17).Value = 141.2708778385117
$ws.Cells.Item(11, 18).Value = 1271.437900546605
$ws.Cells.Item(11, 19).Value = 0.08064266166679575
$ws.Cells.Item(11, 20).Value = 0.08433396977536448

$ws.Cells.Item(12, 7).Value = 54.70735966666666
$ws.Cells.Item(12, 8).Value = 164.122079
$ws.Cells.Item(12, 9).Value = 0.1481428893434501
$ws.Cells.Item(12, 10).Value = 0.1482645698807303
$ws.Cells.Item(12, 13).Value = 1.986769333333333
$ws.Cells.Item(12, 14).Value = 5.960307999999999
$ws.Cells.Item(12, 15).Value = 0.4188172894723593
$ws.Cells.Item(12, 16).Value = 0.4376286247783243
$ws.Cells.Item(12, 17).Value = 108.6909044933702
$ws.Cells.Item(12, 18).Value = 978.2181404403318
$ws.Cells.Item(12, 19).Value = 0.06204480336942744
$ws.Cells.Item(12, 20).Value = 0.06488481982025375

$ws.Cells.Item(13, 7).Value = 54.70735966666666
$ws.Cells.Item(13, 8).Value = 164.122079
$ws.Cells.Item(13, 9).Value = 0.1481428893434501
$ws.Cells.Item(13, 10).Value = 0.1482645698807303
$ws.Cells.Item(13, 15).Value = 0.2760769767479527
$ws.Cells.Item(13, 16).Value = 0.2884770774849725
$ws.Cells.Item(13, 17).Value = 71.64712887171888
$ws.Cells.Item(13, 18).Value = 644.8241598454699
$ws.Cells.Item(13, 19).Value = 0.04089884101664622
$ws.Cells.Item(13, 20).Value = 0.04277092981375955

$ws.Cells.Item(14, 7).Value = 54.70735966666666
$ws.Cells.Item(14, 8).Value = 164.122079
$ws.Cells.Item(14, 9).Value = 0.1481428893434501
$ws.Cells.Item(14, 10).Value = 0.1482645698807303
$ws.Cells.Item(14, 13).Value = 0.3499836666666667
$ws.Cells.Item(14, 14).Value = 1.049951
$ws.Cells.Item(14, 15).Value = 0.07377766919071853
$ws.Cells.Item(14, 16).Value = 0.07709142081493547
$ws.Cells.Item(14, 17).Value = 19.14668232979211
$ws.Cells.Item(14, 18).Value = 172.320140968129
$ws.Cells.Item(14, 19).Value = 0.01092963708293829
$ws.Cells.Item(14, 20).Value = 0.01142992634862078

$ws.Cells.Item(15, 7).Value = 54.70735966666666
$ws.Cells.Item(15, 8).Value = 164.122079
$ws.Cells.Item(15, 9).Value = 0.1481428893434501
$ws.Cells.Item(15, 10).Value = 0.1482645698807303
$ws.Cells.Item(15, 13).Value = 0.6117275
$ws.Cells.Item(15, 14).Value = 1.223455
$ws.Cells.Item(15, 15).Value = 0.1289541010862372
$ws.Cells.Item(15, 16).Value = 0.08983074853315715
$ws.Cells.Item(15, 17).Value = 33.46599636049083
$ws.Cells.Item(15, 18).Value = 200.795978162945
$ws.Cells.Item(15, 19).Value = 0.01910363312760252
$ws.Cells.Item(15, 20).Value = 0.01331871729333258

$ws.Cells.Item(16, 7).Value = 54.70735966666666
$ws.Cells.Item(16, 8).Value = 164.122079
$ws.Cells.Item(16, 9).Value = 0.1481428893434501
$ws.Cells.Item(16, 10).Value = 0.1482645698807303
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 0.4856376666666667
$ws.Cells.Item(16, 14).Value = 1.456913
$ws.Cells.Item(16, 15).Value = 0.1023739635027323
$ws.Cells.Item(16, 16).Value = 0.1069721283886106
$ws.Cells.Item(16, 17).Value = 26.56795449801411
$ws.Cells.Item(16, 18).Value = 239.111590482127
$ws.Cells.Item(16, 19).Value = 0.01516597474683568
$ws.Cells.Item(16, 20).Value = 0.0158601766047636

$ws.Cells.Item(17, 7).Value = 0.909222
$ws.Cells.Item(17, 8).Value = 1.818444
$ws.Cells.Item(17, 9).Value = 0.002462096049879378
$ws.Cells.Item(17, 10).Value = 0.001642745565709015
$ws.Cells.Item(17, 13).Value = 1.986769333333333
$ws.Cells.Item(17, 14).Value = 5.960307999999999
$ws.Cells.Item(17, 15).Value = 0.4188172894723593
$ws.Cells.Item(17, 16).Value = 0.4376286247783243
$ws.Cells.Item(17, 17).Value = 1.806414386792
$ws.Cells.Item(17, 18).Value = 10.838486320752
$ws.Cells.Item(17, 19).Value = 0.001031168394031084
$ws.Cells.Item(17, 20).Value = 0.0007189124827819267

$ws.Cells.Item(18, 7).Value = 0.909222
$ws.Cells.Item(18, 8).Value = 1.818444
$ws.Cells.Item(18, 9).Value = 0.002462096049879378
$ws.Cells.Item(18, 10).Value = 0.001642745565709015
$ws.Cells.Item(18, 15).Value = 0.2760769767479527
$ws.Cells.Item(18, 16).Value = 0.2884770774849725
$ws.Cells.Item(18, 17).Value = 1.19075653082
$ws.Cells.Item(18, 18).Value = 7.144539184919999
$ws.Cells.Item(18, 19).Value = 0.0006797280339137754
$ws.Cells.Item(18, 20).Value = 0.0004738944398471346

$ws.Cells.Item(19, 7).Value = 0.909222
$ws.Cells.Item(19, 8).Value = 1.818444
$ws.Cells.Item(19, 9).Value = 0.002462096049879378
$ws.Cells.Item(19, 10).Value = 0.001642745565709015
$ws.Cells.Item(19, 13).Value = 0.3499836666666667
$ws.Cells.Item(19, 14).Value = 1.049951
$ws.Cells.Item(19, 15).Value = 0.07377766919071853
$ws.Cells.Item(19, 16).Value = 0.07709142081493547
$ws.Cells.Item(19, 17).Value = 0.318212849374
$ws.Cells.Item(19, 18).Value = 1.909277096244
$ws.Cells.Item(19, 19).Value = 0.0001816477078837756
$ws.Cells.Item(19, 20).Value = 0.0001266415896979429

$ws.Cells.Item(20, 7).Value = 0.909222
$ws.Cells.Item(20, 8).Value = 1.818444
$ws.Cells.Item(20, 9).Value = 0.002462096049879378
$ws.Cells.Item(20, 10).Value = 0.001642745565709015
$ws.Cells.Item(20, 13).Value = 0.6117275
$ws.Cells.Item(20, 14).Value = 1.223455
$ws.Cells.Item(20, 15).Value = 0.1289541010862372
$ws.Cells.Item(20, 16).Value = 0.08983074853315715
$ws.Cells.Item(20, 17).Value = 0.5561961010049999
$ws.Cells.Item(20, 18).Value = 2.22478440402
$ws.Cells.Item(20, 19).Value = 0.0003174973829001706
$ws.Cells.Item(20, 20).Value = 0.0001475690638171655

$ws.Cells.Item(21, 7).Value = 0.909222
$ws.Cells.Item(21, 8).Value = 1.818444
$ws.Cells.Item(21, 9).Value = 0.002462096049879378
$ws.Cells.Item(21, 10).Value = 0.001642745565709015
$ws.Cells.Item(21, 11).Value = 2
$ws.Cells.Item(21, 12).Value = 0.6666666666666666
$ws.Cells.Item(21, 13).Value = 0.4856376666666667
$ws.Cells.Item(21, 14).Value = 1.456913
$ws.Cells.Item(21, 15).Value = 0.1023739635027323
$ws.Cells.Item(21, 16).Value = 0.1069721283886106
$ws.Cells.Item(21, 17).Value = 0.441552450562
$ws.Cells.Item(21, 18).Value = 2.649314703372
$ws.Cells.Item(21, 19).Value = 0.0002520545311505729
$ws.Cells.Item(21, 20).Value = 0.0001757279895648455

$ws.Cells.Item(22, 7).Value = 21.33956566666667
$ws.Cells.Item(22, 8).Value = 64.018697
$ws.Cells.Item(22, 9).Value = 0.05778573366465133
$ws.Cells.Item(22, 10).Value = 0.05783319729352075
$ws.Cells.Item(22, 13).Value = 1.986769333333333
$ws.Cells.Item(22, 14).Value = 5.960307999999999
$ws.Cells.Item(22, 15).Value = 0.4188172894723593
$ws.Cells.Item(22, 16).Value = 0.4376286247783243
$ws.Cells.Item(22, 17).Value = 42.39679465318622
$ws.Cells.Item(22, 18).Value = 381.571151878676
$ws.Cells.Item(22, 19).Value = 0.02420166434360093
$ws.Cells.Item(22, 20).Value = 0.02530946259809699

$ws.Cells.Item(23, 7).Value = 21.33956566666667
$ws.Cells.Item(23, 8).Value = 64.018697
$ws.Cells.Item(23, 9).Value = 0.05778573366465133
$ws.Cells.Item(23, 10).Value = 0.05783319729352075
$ws.Cells.Item(23, 15).Value = 0.2760769767479527
$ws.Cells.Item(23, 16).Value = 0.2884770774849725
$ws.Cells.Item(23, 17).Value = 27.94721991157889
$ws.Cells.Item(23, 18).Value = 251.52497920421
$ws.Cells.Item(23, 19).Value = 0.01595331064929933
$ws.Cells.Item(23, 20).Value = 0.01668355173684669

$ws.Cells.Item(24, 7).Value = 21.33956566666667
$ws.Cells.Item(24, 8).Value = 64.018697
$ws.Cells.Item(24, 9).Value = 0.05778573366465133
$ws.Cells.Item(24, 10).Value = 0.05783319729352075
$ws.Cells.Item(24, 13).Value = 0.3499836666666667
$ws.Cells.Item(24, 14).Value = 1.049951
$ws.Cells.Item(24, 15).Value = 0.07377766919071853
$ws.Cells.Item(24, 16).Value = 0.07709142081493547
$ws.Cells.Item(24, 17).Value = 7.468499437094112
$ws.Cells.Item(24, 18).Value = 67.216494933847
$ws.Cells.Item(24, 19).Value = 0.004263296742253612
$ws.Cells.Item(24, 20).Value = 0.004458443349627995

$ws.Cells.Item(25, 7).Value = 21.33956566666667
$ws.Cells.Item(25, 8).Value = 64.018697
$ws.Cells.Item(25, 9).Value = 0.05778573366465133
$ws.Cells.Item(25, 10).Value = 0.05783319729352075
$ws.Cells.Item(25, 13).Value = 0.6117275
$ws.Cells.Item(25, 14).Value = 1.223455
$ws.Cells.Item(25, 15).Value = 0.1289541010862372
$ws.Cells.Item(25, 16).Value = 0.08983074853315715
$ws.Cells.Item(25, 17).Value = 13.05399915635583
$ws.Cells.Item(25, 18).Value = 78.323994938135
$ws.Cells.Item(25, 19).Value = 0.007451707340333825
$ws.Cells.Item(25, 20).Value = 0.005195199402942726

$ws.Cells.Item(26, 7).Value = 21.33956566666667
$ws.Cells.Item(26, 8).Value = 64.018697
$ws.Cells.Item(26, 9).Value = 0.05778573366465133
$ws.Cells.Item(26, 10).Value = 0.05783319729352075
$ws.Cells.Item(26, 11).Value = 2
$ws.Cells.Item(26, 12).Value = 0.6666666666666666
$ws.Cells.Item(26, 13).Value = 0.4856376666666667
$ws.Cells.Item(26, 14).Value = 1.456913
$ws.Cells.Item(26, 15).Value = 0.1023739635027323
$ws.Cells.Item(26, 16).Value = 0.1069721283886106
$ws.Cells.Item(26, 17).Value = 10.36329687804011
$ws.Cells.Item(26, 18).Value = 93.26967190236101
$ws.Cells.Item(26, 19).Value = 0.005915754589163625
$ws.Cells.Item(26, 20).Value = 0.006186540206006348
